# renamed field/value mapping sheets for excel/csv consistency
$wb = $excel.ActiveWorkbook

# 1. Rename the mapping sheets.
$wb.Worksheets.Item("field_mapping").Name = "fields"
$wb.Worksheets.Item("value_mapping").Name = "values"

# 2. Re-write the "data" sheet's D2:D28 column as a single shared formula
#    (D2 becomes the master formula, D3:D28 reuse it) while keeping the
#    existing General-number-format style (s="4") instead of the
#    date-format style Excel would otherwise infer from the C/B operands.
$ws = $wb.Worksheets.Item("data")
$ws.Range("D2:D28").Formula = "=C2-B2"
$ws.Range("D2:D28").NumberFormat = "general"

# 3. Make "values" the active sheet/tab (was "data" before).
$wb.Worksheets.Item("values").Activate()
